$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 61.91334466666666
$ws.Range("H2").Value = 185.740034
$ws.Range("I2").Value = 0.5020829437194911
$ws.Range("J2").Value = 0.5020829437194911
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 91.67302333333333
$ws.Range("N2").Value = 275.01907
$ws.Range("O2").Value = 0.8966431814716052
$ws.Range("P2").Value = 0.896643181471605
$ws.Range("Q2").Value = 5675.783490272041
$ws.Range("R2").Value = 51082.05141244837
$ws.Range("S2").Value = 0.4501892480192734
$ws.Range("T2").Value = 0.4501892480192733

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 61.91334466666666
$ws.Range("H3").Value = 185.740034
$ws.Range("I3").Value = 0.5020829437194911
$ws.Range("J3").Value = 0.5020829437194911
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.548386
$ws.Range("N3").Value = 25.645158
$ws.Range("O3").Value = 0.08361076945850333
$ws.Range("P3").Value = 0.08361076945850332
$ws.Range("Q3").Value = 529.259168761708
$ws.Range("R3").Value = 4763.332518855372
$ws.Range("S3").Value = 0.04197954125637707
$ws.Range("T3").Value = 0.04197954125637707

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 61.91334466666666
$ws.Range("H4").Value = 185.740034
$ws.Range("I4").Value = 0.5020829437194911
$ws.Range("J4").Value = 0.5020829437194911
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.018841
$ws.Range("N4").Value = 6.056523
$ws.Range("O4").Value = 0.01974604906989159
$ws.Range("P4").Value = 0.01974604906989159
$ws.Range("Q4").Value = 124.993198660198
$ws.Range("R4").Value = 1124.938787941782
$ws.Range("S4").Value = 0.009914154443840691
$ws.Range("T4").Value = 0.00991415444384069

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 50.09443433333333
$ws.Range("H5").Value = 150.283303
$ws.Range("I5").Value = 0.4062381250674705
$ws.Range("J5").Value = 0.4062381250674706
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 91.67302333333333
$ws.Range("N5").Value = 275.01907
$ws.Range("O5").Value = 0.8966431814716052
$ws.Range("P5").Value = 0.896643181471605
$ws.Range("Q5").Value = 4592.308247509801
$ws.Range("R5").Value = 41330.7742275882
$ws.Range("S5").Value = 0.3642506448955566
$ws.Range("T5").Value = 0.3642506448955566

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 50.09443433333333
$ws.Range("H6").Value = 150.283303
$ws.Range("I6").Value = 0.4062381250674705
$ws.Range("J6").Value = 0.4062381250674706
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 8.548386
$ws.Range("N6").Value = 25.645158
$ws.Range("O6").Value = 0.08361076945850333
$ws.Range("P6").Value = 0.08361076945850332
$ws.Range("Q6").Value = 428.226561132986
$ws.Range("R6").Value = 3854.039050196874
$ws.Range("S6").Value = 0.03396588222027092
$ws.Range("T6").Value = 0.03396588222027092

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 50.09443433333333
$ws.Range("H7").Value = 150.283303
$ws.Range("I7").Value = 0.4062381250674705
$ws.Range("J7").Value = 0.4062381250674706
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.018841
$ws.Range("N7").Value = 6.056523
$ws.Range("O7").Value = 0.01974604906989159
$ws.Range("P7").Value = 0.01974604906989159
$ws.Range("Q7").Value = 101.132697903941
$ws.Range("R7").Value = 910.1942811354689
$ws.Range("S7").Value = 0.008021597951643032
$ws.Range("T7").Value = 0.008021597951643032

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 11.30520233333333
$ws.Range("H8").Value = 33.915607
$ws.Range("I8").Value = 0.09167893121303822
$ws.Range("J8").Value = 0.09167893121303823
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 91.67302333333333
$ws.Range("N8").Value = 275.01907
$ws.Range("O8").Value = 0.8966431814716052
$ws.Range("P8").Value = 0.896643181471605
$ws.Range("Q8").Value = 1036.382077291721
$ws.Range("R8").Value = 9327.438695625491
$ws.Range("S8").Value = 0.08220328855677503
$ws.Range("T8").Value = 0.08220328855677503

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 11.30520233333333
$ws.Range("H9").Value = 33.915607
$ws.Range("I9").Value = 0.09167893121303822
$ws.Range("J9").Value = 0.09167893121303823
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 8.548386
$ws.Range("N9").Value = 25.645158
$ws.Range("O9").Value = 0.08361076945850333
$ws.Range("P9").Value = 0.08361076945850332
$ws.Range("Q9").Value = 96.64123335343402
$ws.Range("R9").Value = 869.7711001809062
$ws.Range("S9").Value = 0.007665345981855324
$ws.Range("T9").Value = 0.007665345981855324

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 11.30520233333333
$ws.Range("H10").Value = 33.915607
$ws.Range("I10").Value = 0.09167893121303822
$ws.Range("J10").Value = 0.09167893121303823
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.018841
$ws.Range("N10").Value = 6.056523
$ws.Range("O10").Value = 0.01974604906989159
$ws.Range("P10").Value = 0.01974604906989159
$ws.Range("Q10").Value = 22.823405983829
$ws.Range("R10").Value = 205.410653854461
$ws.Range("S10").Value = 0.001810296674407869
$ws.Range("T10").Value = 0.001810296674407869
